$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property")

$ws.Range("A31").Value = "DefaultHitTime"
$ws.Range("B31").Value = "float"
$ws.Range("C31").Value = $false
$ws.Range("D31").Value = $false
$ws.Range("E31").Value = $false
$ws.Range("F31").Value = $true
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = "Friend"
$ws.Range("J31").Value = "缺省打击时间（本来应该打到但是物理没碰撞到或者其他原因）"

# Match the text-formatted style used by the other label/desc columns (A, B, I, J)
$ws.Range("A31").NumberFormat = "@"
$ws.Range("B31").NumberFormat = "@"
$ws.Range("I31").NumberFormat = "@"
$ws.Range("J31").NumberFormat = "@"
